# More Ellipse intersection fixes
# Fixing some cases where intersections between oblique ellipses were failing.
#
# Renames "Rotated Ellipse" -> "Oblique Ellipse" and
#         "Rotated Elliptical Arc" -> "Oblique Elliptical Arc"
# throughout the "Intersection" sheet/table (header row + row labels),
# and flips a few intersection results that were fixed by the underlying
# code change (Quadratic/Cubic Bezier vs Oblique Ellipse, now both TRUE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intersection")

# Header row (row 1) - rename the two "Rotated" columns to "Oblique"
$ws.Range("O1").Value = "Oblique Ellipse"
$ws.Range("Q1").Value = "Oblique Elliptical Arc"

# Row labels (column A) - same rename, mirrored on the transposed axis
$ws.Range("A15").Value = "Oblique Ellipse"
$ws.Range("A17").Value = "Oblique Elliptical Arc"

# Intersection matrix fixes: Quadratic Bezier / Cubic Bezier now
# correctly intersect with the Oblique Ellipse (both directions of the
# symmetric matrix get updated).
$ws.Range("O6").Value = 1
$ws.Range("O7").Value = 1
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1

# Leave the cursor on A18, matching the saved selection state.
$ws.Range("A18").Select() | Out-Null
